# Generate Report for Handback
# -----------------------------------------------------------------------
# The file "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md" has now also been
# handed back (status "Handed back: in sync with en-US") alongside
# "9e9b3e7d-d93e-447a-bbad-150e428577a6.md". Both the zh-cn and de-de
# per-language sheets as well as the Overview roll-up sheet are updated
# so that the 3b5cfbe7 row leads (row 2) and the 9e9b3e7d row follows
# (row 3), matching the new handback ordering, and the new handback
# timestamps are recorded.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("D2").Value = "2016-03-24 08:27:02"

$overview.Range("A3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"
$overview.Range("D3").Value = "2016-03-24 08:25:29"

foreach ($hl in $overview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
    }
}

# --- Sheet "zh-cn" --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
$zhcn.Range("B2").Value = ".md"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("D2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf"
$zhcn.Range("E2").Value = "2016-03-24 08:26:57"
$zhcn.Range("F2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
$zhcn.Range("G2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-03-24 08:27:22"
$zhcn.Range("J2").Value = "Include"

$zhcn.Range("A3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("D3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-24 08:25:24"
$zhcn.Range("F3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
$zhcn.Range("G3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-03-24 08:26:03"
$zhcn.Range("J3").Value = "Include"

foreach ($hl in $zhcn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf"
    } elseif ($addr -eq '$F$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
    } elseif ($addr -eq '$G$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf"
    } elseif ($addr -eq '$F$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
    } elseif ($addr -eq '$G$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf"
    }
}

# --- Sheet "de-de" --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
$dede.Range("B2").Value = ".md"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("D2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf"
$dede.Range("E2").Value = "2016-03-24 08:27:02"
$dede.Range("F2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
$dede.Range("G2").Value = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf"
$dede.Range("H2").Value = "2016-03-24 08:27:29"
$dede.Range("J2").Value = "Include"

$dede.Range("A3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("D3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf"
$dede.Range("E3").Value = "2016-03-24 08:25:29"
$dede.Range("F3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
$dede.Range("G3").Value = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf"
$dede.Range("H3").Value = "2016-03-24 08:26:14"
$dede.Range("J3").Value = "Include"

foreach ($hl in $dede.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf"
    } elseif ($addr -eq '$F$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md"
    } elseif ($addr -eq '$G$2') {
        $hl.TextToDisplay = "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf"
    } elseif ($addr -eq '$F$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.md"
    } elseif ($addr -eq '$G$3') {
        $hl.TextToDisplay = "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf"
    }
}

Write-Output "Report generated for handback."
